# Update the "label" column (M) values for the three variable rows on
# Sheet1 so they hold the literal column header text instead of stale
# per-row captions (catches non-initialised item access before the
# item is prepared for access).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Range("M3").Value = "label"
$ws1.Range("M4").Value = "label"
$ws1.Range("M5").Value = "label"

# Leave Sheet2's selection where it was, then select/activate Sheet1 at
# G7 last so Sheet1 ends up the active tab.
$ws2.Range("A2").Select()
$ws1.Range("G7").Select()
